$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-observation row needs to be inserted at row 300 (just above
# the existing row 300), pushing all subsequent rows (300-353) down by one
# (new rows 301-354). The new row carries the latest "Ajo" (garlic) price
# reading; the rest of the data below is unchanged, just shifted down.

$ws.Rows.Item(300).Insert()

$ws.Cells.Item(300, 1).Value = 6
$ws.Cells.Item(300, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(300, 3).Value = "Metropolitana"
$ws.Cells.Item(300, 4).Value = 44476
$ws.Cells.Item(300, 5).Value = 13
$ws.Cells.Item(300, 6).Value = 100112003
$ws.Cells.Item(300, 7).Value = "Ajo"
$ws.Cells.Item(300, 8).Value = "Chino"
$ws.Cells.Item(300, 9).Value = "Primera"
$ws.Cells.Item(300, 10).Value = 2300
$ws.Cells.Item(300, 11).Value = 14500
$ws.Cells.Item(300, 12).Value = 15000
$ws.Cells.Item(300, 13).Value = 14783
$ws.Cells.Item(300, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(300, 15).Value = "China"
$ws.Cells.Item(300, 16).Value = 1478
$ws.Cells.Item(300, 17).Value = 10
$ws.Cells.Item(300, 18).Value = "Hortaliza"
